$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for specific rows, per repulled data.
$ws.Range("F10").Value = -3
$ws.Range("F12").Value = -2
$ws.Range("F13").Value = -3
$ws.Range("F14").Value = -1
$ws.Range("F21").Value = 0
$ws.Range("F23").Value = -2
$ws.Range("F27").Value = -1
$ws.Range("F29").Value = 4
$ws.Range("F32").Value = -2
$ws.Range("F34").Value = -1
$ws.Range("F38").Value = 1
$ws.Range("F40").Value = 5
$ws.Range("F43").Value = 0
